$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 333336260
$ws.Range("I18").Value = 3890
$ws.Range("K18").Value = 3890
$ws.Range("M18").Value = -3606

$ws.Range("H40").Value = 5190.2
$ws.Range("I40").Value = 4916.5
$ws.Range("J40").Value = 5600.75
$ws.Range("K40").Value = 4916.5
$ws.Range("L40").Value = 5600.75
$ws.Range("M40").Value = -4741.5
$ws.Range("N40").Value = -5950.75

$ws.Range("H107").Value = 796.06665
$ws.Range("I107").Value = 726.4167
$ws.Range("J107").Value = 1074.6666
$ws.Range("K107").Value = 726.4167
$ws.Range("L107").Value = 1074.6666
$ws.Range("M107").Value = 1193.5833
$ws.Range("N107").Value = -4914.6666

$ws.Range("H112").Value = 1527.2916
$ws.Range("J112").Value = 1572.6086
$ws.Range("L112").Value = 4717.825800000001
$ws.Range("N112").Value = -6933.825800000001

$ws.Range("H125").Value = 927.9
$ws.Range("I125").Value = 931.125
$ws.Range("J125").Value = 915
$ws.Range("K125").Value = 8380.125
$ws.Range("L125").Value = 8235
$ws.Range("M125").Value = -5920.125
$ws.Range("N125").Value = -13155

$ws.Range("H132").Value = 7319.8945
$ws.Range("I132").Value = 7615.1113
$ws.Range("K132").Value = 22845.3339
$ws.Range("M132").Value = -20315.3339

$ws.Range("H137").Value = 2268.6
$ws.Range("I137").Value = 2230.5186
$ws.Range("J137").Value = 2397.125
$ws.Range("K137").Value = 6691.5558
$ws.Range("L137").Value = 7191.375
$ws.Range("M137").Value = -4141.5558
$ws.Range("N137").Value = -12291.375

$ws.Range("H141").Value = 2885.923
$ws.Range("I141").Value = 2626.4167
$ws.Range("K141").Value = 7879.250100000001
$ws.Range("M141").Value = -2699.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1926443.2
$ws.Range("I32").Value = 863824.5600000001
$ws.Range("K32").Value = 863824.5600000001
$ws.Range("M32").Value = -863537.5600000001

$ws.Range("H61").Value = 3557
$ws.Range("I61").Value = 2671.8
$ws.Range("J61").Value = 3999.6
$ws.Range("K61").Value = 2671.8
$ws.Range("L61").Value = 3999.6
$ws.Range("M61").Value = -2459.8
$ws.Range("N61").Value = -4423.6

$ws.Range("H74").Value = 1893.8334
$ws.Range("I74").Value = 1455.4706
$ws.Range("J74").Value = 2958.4285
$ws.Range("K74").Value = 1455.4706
$ws.Range("L74").Value = 2958.4285
$ws.Range("M74").Value = -581.4706000000001
$ws.Range("N74").Value = -4706.4285

$ws.Range("H77").Value = 1893.8334
$ws.Range("I77").Value = 1455.4706
$ws.Range("J77").Value = 2958.4285
$ws.Range("K77").Value = 7277.353000000001
$ws.Range("L77").Value = 14792.1425
$ws.Range("M77").Value = -2909.353000000001
$ws.Range("N77").Value = -23528.1425

$ws.Range("H122").Value = 3056.5334
$ws.Range("I122").Value = 3196.077
$ws.Range("J122").Value = 2149.5
$ws.Range("K122").Value = 9588.231
$ws.Range("L122").Value = 6448.5
$ws.Range("M122").Value = -7138.231
$ws.Range("N122").Value = -11348.5

$ws.Range("H132").Value = 4752.7646
$ws.Range("I132").Value = 8559.6
$ws.Range("J132").Value = 3166.5833
$ws.Range("K132").Value = 25678.8
$ws.Range("L132").Value = 9499.749899999999
$ws.Range("M132").Value = -23148.8
$ws.Range("N132").Value = -14559.7499

$ws.Range("H136").Value = 3557
$ws.Range("I136").Value = 2671.8
$ws.Range("J136").Value = 3999.6
$ws.Range("K136").Value = 8015.400000000001
$ws.Range("L136").Value = 11998.8
$ws.Range("M136").Value = -5465.400000000001
$ws.Range("N136").Value = -17098.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2175.652
$ws.Range("I134").Value = 1341.25
$ws.Range("K134").Value = 4023.75
$ws.Range("M134").Value = -1488.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 26633.334
$ws.Range("I41").Value = 26633.334
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 26633.334
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -26205.334
$ws.Range("N41").Value = $null

$ws.Range("H50").Value = 59666
$ws.Range("J50").Value = 59666
$ws.Range("L50").Value = 59666
$ws.Range("N50").Value = -60916

$ws.Range("H60").Value = 23819.6
$ws.Range("J60").Value = 24774.5
$ws.Range("L60").Value = 24774.5
$ws.Range("N60").Value = -25796.5

$ws.Range("H99").Value = 5998.5
$ws.Range("J99").Value = 6498.75
$ws.Range("L99").Value = 6498.75
$ws.Range("N99").Value = -9494.75

$ws.Range("H105").Value = 2875
$ws.Range("I105").Value = 2805.5557
$ws.Range("K105").Value = 2805.5557
$ws.Range("M105").Value = -1058.5557

$ws.Range("H126").Value = 5998.5
$ws.Range("J126").Value = 6498.75
$ws.Range("L126").Value = 19496.25
$ws.Range("N126").Value = -24436.25

$ws.Range("H132").Value = 4121.5
$ws.Range("I132").Value = 3579.4
$ws.Range("J132").Value = 5025
$ws.Range("K132").Value = 10738.2
$ws.Range("L132").Value = 15075
$ws.Range("M132").Value = -8208.200000000001
$ws.Range("N132").Value = -20135

$ws.Range("H134").Value = 3760.9062
$ws.Range("I134").Value = 4204.2915
$ws.Range("K134").Value = 12612.8745
$ws.Range("M134").Value = -10077.8745

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1903.5
$ws.Range("J113").Value = 1903.5
$ws.Range("L113").Value = 5710.5
$ws.Range("N113").Value = -10050.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5384.1665
$ws.Range("I70").Value = 4557
$ws.Range("K70").Value = 4557
$ws.Range("M70").Value = -4287

$ws.Range("H73").Value = 5384.1665
$ws.Range("I73").Value = 4557
$ws.Range("K73").Value = 4557
$ws.Range("M73").Value = -3621

$ws.Range("H80").Value = 111115160
$ws.Range("I80").Value = 200002300
$ws.Range("K80").Value = 200002300
$ws.Range("M80").Value = -200001302

$ws.Range("H83").Value = 111115160
$ws.Range("I83").Value = 200002300
$ws.Range("K83").Value = 1000011500
$ws.Range("M83").Value = -1000006508

$ws.Range("H102").Value = 1394.4103
$ws.Range("I102").Value = 1125.8667
$ws.Range("J102").Value = 1562.25
$ws.Range("K102").Value = 1125.8667
$ws.Range("L102").Value = 1562.25
$ws.Range("M102").Value = 496.1333
$ws.Range("N102").Value = -4806.25

$ws.Range("H107").Value = 5840.5713
$ws.Range("J107").Value = 8998
$ws.Range("L107").Value = 8998
$ws.Range("N107").Value = -12838

$ws.Range("H132").Value = 2112.389
$ws.Range("I132").Value = 1447.5714
$ws.Range("K132").Value = 4342.7142
$ws.Range("M132").Value = -1812.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1123.8
$ws.Range("I22").Value = 783.3333
$ws.Range("J22").Value = 1634.5
$ws.Range("K22").Value = 783.3333
$ws.Range("L22").Value = 1634.5
$ws.Range("M22").Value = -488.3333
$ws.Range("N22").Value = -2224.5

$ws.Range("H27").Value = 1123.8
$ws.Range("I27").Value = 783.3333
$ws.Range("J27").Value = 1634.5
$ws.Range("K27").Value = 783.3333
$ws.Range("L27").Value = 1634.5
$ws.Range("M27").Value = -676.3333
$ws.Range("N27").Value = -1848.5

$ws.Range("H40").Value = 16901.25
$ws.Range("I40").Value = 20666.666
$ws.Range("K40").Value = 20666.666
$ws.Range("M40").Value = -20530.666

$ws.Range("H46").Value = 2287.7693
$ws.Range("I46").Value = 2018.375
$ws.Range("K46").Value = 2018.375
$ws.Range("M46").Value = -1830.375

$ws.Range("H122").Value = 4696.852
$ws.Range("I122").Value = 3350.6875
$ws.Range("K122").Value = 10052.0625
$ws.Range("M122").Value = -7602.0625

$ws.Range("H132").Value = 4016.3635
$ws.Range("I132").Value = 4910
$ws.Range("K132").Value = 14730
$ws.Range("M132").Value = -12200

$ws.Range("H136").Value = 7695.2
$ws.Range("J136").Value = 7079
$ws.Range("L136").Value = 21237
$ws.Range("N136").Value = -26337

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 33165.668
$ws.Range("I39").Value = 3500
$ws.Range("J39").Value = 47998.5
$ws.Range("K39").Value = 3500
$ws.Range("L39").Value = 47998.5
$ws.Range("M39").Value = -3087
$ws.Range("N39").Value = -48824.5

$ws.Range("H122").Value = 13890781
$ws.Range("I122").Value = 2047.6666
$ws.Range("J122").Value = 41668250
$ws.Range("K122").Value = 6142.9998
$ws.Range("L122").Value = 125004750
$ws.Range("M122").Value = -3692.9998
$ws.Range("N122").Value = -125009650
